$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F3..F6
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 225
$wsExhibit.Range("F4").Value = 841
$wsExhibit.Range("F5").Value = 73
$wsExhibit.Range("F6").Value = 30

# Sheet "全部类型" (sheet4): F4..F7
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 225
$wsAll.Range("F5").Value = 841
$wsAll.Range("F6").Value = 73
$wsAll.Range("F7").Value = 30
